# Apply the CELGNamed.xlsx update:
#  - Add two new data points to row 3 (X3, Y3) classifying the move as "Down"
#  - Append a brand new row 4 of sentiment/trading data (a "Buy"/"Named"/"N/A" record)
#  - The "Down" string becomes a brand-new shared string entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: append the new percent-change / direction columns ---
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = -0.93999500000001035
$ws.Range("Y3").Value = "Down"

# --- Row 4: brand new record appended beneath row 3 ---
$ws.Range("A4").Value = 42633.890567129631
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 44
$ws.Range("E4").Value = 14226
$ws.Range("F4").Value = 807
$ws.Range("G4").Value = 67
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 16504
$ws.Range("L4").Value = 154
$ws.Range("M4").Value = 70
$ws.Range("N4").Value = 16
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Named"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.86
$ws.Range("S4").Value = 0.0262
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = -2.66
$ws.Range("U4").Value = 15.05
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
